# Update crypto price/volume snapshot data (GitHub Actions scrape refresh).
# Two coin rows (17/18 and 49/50) were also re-ranked and swapped in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.230.57'
$ws.Range("E2").Value = '  +2.59%  '
# Row 3
$ws.Range("D3").Value = '2.595.58'
$ws.Range("E3").Value = '  +1.64%  '
# Row 4
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.19%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.40%  '
# Row 7
$ws.Range("E7").Value = '  +0.00%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.567'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.69%  '
# Row 9
$ws.Range("D9").Value = '2.609.90'
$ws.Range("E9").Value = '  +1.47%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '
# Row 11
$ws.Range("E11").Value = '  +4.54%  '
# Row 12
$ws.Range("E12").Value = '  +3.89%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.135'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.66%  '
# Row 14
$ws.Range("D14").Value = '3.055.74'
$ws.Range("E14").Value = '  +1.78%  '
# Row 15
$ws.Range("D15").Value = '59.180.03'
$ws.Range("E15").Value = '  +2.52%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.21%  '
# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.597.85'
$ws.Range("E17").Value = '  +1.54%  '
# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.65%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '346.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.75%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.12%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.28%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.11%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.82%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.03%  '
# Row 26
$ws.Range("E26").Value = '  +3.11%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.69%  '
# Row 29
$ws.Range("E29").Value = '  +0.06%  '
# Row 30
$ws.Range("D30").Value = '0.0₃0735'
$ws.Range("E30").Value = '  +5.40%  '
# Row 31
$ws.Range("E31").Value = '  +5.60%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.92%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.80'
$ws.Range("D33").Style = "Normal"
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.49%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.54%  '
# Row 36
$ws.Range("E36").Value = '  +2.06%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.95'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.77%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.47'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.78%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.837'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.34%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.828'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '
# Row 41
$ws.Range("E41").Value = '  +2.50%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '276.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.86%  '
# Row 43
$ws.Range("E43").Value = '  -0.04%  '
# Row 44
$ws.Range("E44").Value = '  +3.05%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.59%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0961'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.97%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0522'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.52%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.07%  '
# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0223'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.09%  '
# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.940.48'
$ws.Range("E50").Value = '  -0.17%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.02%  '
